$wb = $excel.ActiveWorkbook

# --- Sheet "Overview" ---
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E3").Value = "Ready for handoff"
$overview.Range("F3").Value = "Ready for handoff"
$overview.Range("G3").Value = "2016-08-29 16:53:21"

# --- Sheet "zh-cn" ---
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C2").Value = "Ready for handoff"
$zhcn.Range("C3").Value = "Ready for handoff"
$zhcn.Range("H3").Value = "2016-08-29 16:53:16"
$zhcn.Range("P3").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/fc965fade1f0e77c461275c0f99463fe41e21995/e2e/fd31c7fa-9349-463d-91cb-649c56cef66f.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/6c8a555fcf4f1baa5c447b7f8679386b35356a7e/e2e/fd31c7fa-9349-463d-91cb-649c56cef66f.md."
$zhcn.Columns.Item(16).ColumnWidth = 39.17

# --- Sheet "de-de" ---
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C2").Value = "Ready for handoff"
$dede.Range("C3").Value = "Ready for handoff"
$dede.Range("H3").Value = "2016-08-29 16:53:21"
$dede.Range("P3").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/fc965fade1f0e77c461275c0f99463fe41e21995/e2e/fd31c7fa-9349-463d-91cb-649c56cef66f.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/6c8a555fcf4f1baa5c447b7f8679386b35356a7e/e2e/fd31c7fa-9349-463d-91cb-649c56cef66f.md."
$dede.Columns.Item(16).ColumnWidth = 39.17
